# CIERRE 9 DIC 2021
# Update the "REMISIONES NOVIEMBRE 2021" credit-tracking sheet with the
# payments / collections recorded since the last close, and move the
# active selection to reflect where the user left off (A21).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("REMISIONES NOVIEMBRE   2021   ")

# Dates are written as raw serial numbers (not [DateTime]) so Excel does
# not rewrite the cell's existing date number-format style.

# Row 15 - GUSTAVO 2500 collected on 28-Nov-2021
$ws.Range("F15").Value = 44528
$ws.Range("G15").Value = 2500

# Row 16 - GUSTAVO 6243 collected on 26-Nov-2021
$ws.Range("F16").Value = 44526
$ws.Range("G16").Value = 6243

# Row 17 - new entry: DAVID HERRADURA, remitted 26-Nov, 6699, paid 27-Nov
$ws.Range("A17").Value = 44526
$ws.Range("D17").Value = "DAVID HERRADURA"
$ws.Range("E17").Value = 6699
$ws.Range("F17").Value = 44527
$ws.Range("G17").Value = 6699

# Row 18 - new entry: OBRADOR, remitted 27-Nov, 351, paid 27-Nov
$ws.Range("A18").Value = 44527
$ws.Range("D18").Value = "OBRADOR"
$ws.Range("E18").Value = 351
$ws.Range("F18").Value = 44527
$ws.Range("G18").Value = 351

# Row 19 - new entry: GUSTAVO, remitted 28-Nov, 5312, still unpaid
$ws.Range("A19").Value = 44528
$ws.Range("D19").Value = "GUSTAVO"
$ws.Range("E19").Value = 5312

# Row 20 - new entry: MAURO, remitted 28-Nov, 975, paid 29-Nov
$ws.Range("A20").Value = 44528
$ws.Range("D20").Value = "MAURO"
$ws.Range("E20").Value = 975
$ws.Range("F20").Value = 44529
$ws.Range("G20").Value = 975

# Reflect the scrolled/selected state left by the user after closing out.
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("A21").Select()
